# ================================================================
# Kpi_General.xlsx update
#   - Switch base fonts Calibri -> Arial
#   - Add new "Times New Roman" styled rows on "Chi tieu" sheet
#   - Refresh the KPI template placeholders on "KPI nhan vien" sheet
#   - Update the "Chi tieu" indicator list with the new KPI names
# ================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("KPI nhan vien")
$ws2 = $wb.Worksheets.Item("Chi tieu")
$ws3 = $wb.Worksheets.Item("Quy tac import")

# ----------------------------------------------------------------
# 1. Global font swap: Calibri -> Arial everywhere it is currently used
# ----------------------------------------------------------------
$ws1.Cells.Font.Name = "Arial"
$ws2.Cells.Font.Name = "Arial"

# ----------------------------------------------------------------
# 2. "Chi tieu" sheet: replace the indicator list (rows 2-11) with the
#    new KPI names, and blank out the now-unused row 12.
# ----------------------------------------------------------------
$ws2.Cells.Item(2, 1).Value  = "Tổng doanh thu đơn hàng"
$ws2.Cells.Item(3, 1).Value  = "Doanh thu C2 Trọng điểm"
$ws2.Cells.Item(4, 1).Value  = "Doanh thu C2 Siêu lớn"
$ws2.Cells.Item(5, 1).Value  = "Doanh thu C2"
$ws2.Cells.Item(6, 1).Value  = "Tổng số đại lý mở mới"
$ws2.Cells.Item(7, 1).Value  = "Số đại lý trọng điểm mở mới"
$ws2.Cells.Item(8, 1).Value  = "Số đại lý ghé thăm"
$ws2.Cells.Item(9, 1).Value  = "Tổng số lượt ghé thăm"
$ws2.Cells.Item(10, 1).Value = "Số thông tin phản ánh"
$ws2.Cells.Item(11, 1).Value = "Số hình ảnh chụp"
$ws2.Cells.Item(12, 1).ClearContents()

# New font treatment for the indicator list: Times New Roman.
$r2to4 = $ws2.Range($ws2.Cells.Item(2, 1), $ws2.Cells.Item(4, 1))
$r2to4.Font.Name = "Times New Roman"
$r2to4.Font.Size = 11
$r2to4.Font.Family = 1

$r3to4 = $ws2.Range($ws2.Cells.Item(3, 1), $ws2.Cells.Item(4, 1))
$r3to4.WrapText = $true

$r5to11 = $ws2.Range($ws2.Cells.Item(5, 1), $ws2.Cells.Item(11, 1))
$r5to11.Font.Name = "Times New Roman"
$r5to11.Font.Size = 11
$r5to11.Font.Family = 1
$r5to11.Font.Color = 0
$r5to11.HorizontalAlignment = -4131
$r5to11.VerticalAlignment = -4108

# ----------------------------------------------------------------
# 3. "KPI nhan vien" sheet: update the per-indicator template column (C)
#    rows 7-16 so the placeholders point at the new KPI fields.
# ----------------------------------------------------------------
$ws1.Cells.Item(7, 3).Value  = "{{KpiGenerals.TotalIndirectSalesAmount.Name}}"
$ws1.Cells.Item(8, 3).Value  = "{{KpiGenerals.RevenueC2TD.Name}}"
$ws1.Cells.Item(9, 3).Value  = "{{KpiGenerals.RevenueC2SL.Name}}"
$ws1.Cells.Item(10, 3).Value = "{{KpiGenerals.RevenueC2.Name}}"
$ws1.Cells.Item(11, 3).Value = "{{KpiGenerals.NewStoresCreated.Name}}"
$ws1.Cells.Item(12, 3).Value = "{{KpiGenerals.NewStoreC2Created.Name}}"
$ws1.Cells.Item(13, 3).Value = "{{KpiGenerals.StoresVisited.Name}}"
$ws1.Cells.Item(14, 3).Value = "{{KpiGenerals.NumberOfStoreVisits.Name}}"
$ws1.Cells.Item(15, 3).Value = "{{KpiGenerals.TotalProblem.Name}}"
$ws1.Cells.Item(16, 3).Value = "{{KpiGenerals.TotalImage.Name}}"

# Row 17 (the old "SKUDirectOrder" indicator row) no longer exists;
# deleting it shifts the trailing "END" row up from 18 to 17.
$ws1.Rows.Item(17).Delete()

Write-Host "Kpi_General.xlsx updated"
